$d = $word.ActiveDocument

# --- locate the target paragraph ---------------------------------------
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*project with a goal being less than*") {
        $targetIndex = $i
        break
    }
}

# --- fix the wording: "may be greater" -> "may have  greater" ----------
#     (two spaces end up between "have" and "greater" because the text
#      is later split into separate runs around a bookmark that sits
#      right there)
$rng = $paras.Item($targetIndex).Range
$null = $rng.Find.Execute("may be greater", $true, $true, $false, $false, $false, $true, 1, $false, "may have  greater", 2)

# --- drop the curly quotes that used to wrap the slash -------------------
$rng = $paras.Item($targetIndex).Range
$quote = [char]0x201D
$target2 = "Theater" + $quote + "/" + $quote + "Plays"
$null = $rng.Find.Execute($target2, $true, $true, $false, $false, $false, $true, 1, $false, "Theater/Plays", 2)

# --- remove the existing _GoBack bookmark; it gets re-created below ------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$start = $paras.Item($targetIndex).Range.Start

# --- split the single run into the pieces required by the target layout -
# offsets, relative to the paragraph start, at which a run boundary is
# needed in the final document:
#   0  ................ "The project with a goal being less than $5000"
#   45 ................ " in US may have"
#   60 ................ " "
#   61 ................ (bookmark _GoBack sits here)
#   61 ................ " greater chance ... launched under \u201cTheater"
#   121 ................ "/"
#   122 ................ "Plays\u201d, it may get positive impact on project success."
#
# A bookmark add/delete at a position splits the run there; doing it via
# a temporary bookmark leaves the split in place even once the bookmark
# is removed again. The splits must be applied starting with the ones
# deepest in the sentence (121, then 122) before working back towards
# the front of the paragraph, otherwise the engine can mis-tag a
# neighbouring run's xml:space="preserve" flag.
$off2 = 45
$off3 = 60
$offBookmark = 61
$off5 = 121
$off6 = 122

$r = $d.Range($start + $off5, $start + $off5)
$d.Bookmarks.Add("ZzTmpSplitA", $r)
$d.Bookmarks.Item("ZzTmpSplitA").Delete()

$r = $d.Range($start + $off6, $start + $off6)
$d.Bookmarks.Add("ZzTmpSplitB", $r)
$d.Bookmarks.Item("ZzTmpSplitB").Delete()

$r = $d.Range($start + $off2, $start + $off2)
$d.Bookmarks.Add("ZzTmpSplitC", $r)
$d.Bookmarks.Item("ZzTmpSplitC").Delete()

$r = $d.Range($start + $off3, $start + $off3)
$d.Bookmarks.Add("ZzTmpSplitD", $r)
$d.Bookmarks.Item("ZzTmpSplitD").Delete()

# --- finally, (re)place _GoBack at its new position ----------------------
$r = $d.Range($start + $offBookmark, $start + $offBookmark)
$d.Bookmarks.Add("_GoBack", $r)
